$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.476.82'
$ws.Range("E2").Value = '  +1.80%  '
$ws.Range("D3").Value = '1.858.77'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'" + '245.14'
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").Value = "'" + '0.6949'
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'" + '0.07694'
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").Value = "'" + '23.63'
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("D11").Value = "'" + '0.07768'
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").Value = "'" + '5.153'
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("D13").Value = '1.856.69'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = "'" + '0.6924'
$ws.Range("E15").Value = '  +1.66%  '
$ws.Range("D16").Value = "'" + '6.340'
$ws.Range("E16").Value = '  -1.76%  '
$ws.Range("D17").Value = '29.459.49'
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = "'" + '0.000008298'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = '2.100.48'
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").Value = "'" + '238.34'
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("D21").Value = "'" + '12.73'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = "'" + '1.000'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = "'" + '7.642'
$ws.Range("E23").Value = '  +2.05%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("D26").Value = "'" + '8.905'
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("D27").Value = "'" + '159.73'
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").Value = "'" + '1.534'
$ws.Range("E29").Value = '  -1.20%  '
$ws.Range("D30").Value = "'" + '4.252'
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("D32").Value = "'" + '1.215'
$ws.Range("E32").Value = '  +3.70%  '
$ws.Range("D33").Value = "'" + '0.05115'
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("D35").Value = "'" + '1.887'
$ws.Range("D36").Value = "'" + '1.149'
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").Value = "'" + '2.682'
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = '1.331.69'
$ws.Range("E38").Value = '  +7.20%  '
$ws.Range("D39").Value = "'" + '0.01870'
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("E40").Value = '  +0.67%  '
$ws.Range("D41").Value = "'" + '0.9574'
$ws.Range("E41").Value = '  +1.71%  '
$ws.Range("D42").Value = "'" + '5.856'
$ws.Range("E42").Value = '  +2.90%  '
$ws.Range("D43").Value = "'" + '105.84'
$ws.Range("E43").Value = '  -2.35%  '
$ws.Range("D45").Value = "'" + '9.908'
$ws.Range("E45").Value = '  +3.44%  '
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("D47").Value = '1.999.04'
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").Value = "'" + '0.5226'
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("D49").Value = "'" + '1.785'
$ws.Range("E49").Value = '  +2.04%  '
$ws.Range("D50").Value = "'" + '63.38'
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("D51").Value = "'" + '6.976'
$ws.Range("E51").Value = '  +0.72%  '
